# Updates crypto price/volume data per upstream diff (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new text would otherwise be auto-parsed as a number
# by Excel (losing the original plain-text representation, e.g. trailing
# zeros) are first switched to the builtin Text format ("@") so the
# assigned string is preserved exactly, matching the source inlineStr cells.

$ws.Range('D2').Value = '68.226.83'
$ws.Range('E2').Value = '  +2.07%  '

$ws.Range('D3').Value = '3.921.19'
$ws.Range('E3').Value = '  +0.90%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '486.00'
$ws.Range('E5').Value = '  +4.17%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.33'
$ws.Range('E6').Value = '  +4.13%  '

$ws.Range('E7').Value = '  +1.39%  '

$ws.Range('E8').Value = '  -0.10%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.732'
$ws.Range('E9').Value = '  +0.00%  '

$ws.Range('E10').Value = '  +5.40%  '

$ws.Range('E11').Value = '  +7.23%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '42.96'
$ws.Range('E12').Value = '  +0.38%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.63'
$ws.Range('E13').Value = '  +3.42%  '

$ws.Range('D14').Value = '4.546.53'
$ws.Range('E14').Value = '  +0.89%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.76'
$ws.Range('E15').Value = '  -0.28%  '

$ws.Range('D16').Value = '3.933.96'
$ws.Range('E16').Value = '  +1.58%  '

$ws.Range('E17').Value = '  -0.12%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.00'
$ws.Range('E18').Value = '  +1.31%  '

$ws.Range('E19').Value = '  -1.64%  '

$ws.Range('D20').Value = '68.348.39'
$ws.Range('E20').Value = '  +1.94%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '443.51'
$ws.Range('E21').Value = '  +3.47%  '

$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.80'
$ws.Range('E22').Value = '  +1.77%  '

$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.40'
$ws.Range('E23').Value = '  +2.13%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '88.51'
$ws.Range('E24').Value = '  +0.66%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.56'
$ws.Range('E25').Value = '  +15.99%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.11'
$ws.Range('E26').Value = '  +15.19%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.63'
$ws.Range('E27').Value = '  +3.00%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '38.90'
$ws.Range('E28').Value = '  +1.73%  '

$ws.Range('E29').Value = '  +1.96%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '721.11'
$ws.Range('E30').Value = '  -0.83%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.57'
$ws.Range('E31').Value = '  -0.35%  '

$ws.Range('E32').Value = '  +0.51%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.88'
$ws.Range('E33').Value = '  +3.06%  '

$ws.Range('D34').Value = '0.0₃0915'
$ws.Range('E34').Value = '  +17.18%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '41.54'
$ws.Range('E35').Value = '  -3.56%  '

$ws.Range('E36').Value = '  +11.25%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '59.38'
$ws.Range('E37').Value = '  +3.76%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.150'
$ws.Range('E38').Value = '  -3.72%  '

$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.397'
$ws.Range('E39').Value = '  +17.70%  '

$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.26%  '

$ws.Range('E41').Value = '  +15.20%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0482'
$ws.Range('E42').Value = '  +1.72%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.16'
$ws.Range('E43').Value = '  +2.61%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.90'
$ws.Range('E44').Value = '  +3.84%  '

$ws.Range('E45').Value = '  +1.86%  '

$ws.Range('E46').Value = '  -0.03%  '

$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0354'
$ws.Range('E47').Value = '  +45.51%  '

$ws.Range('B48').Value = 'LidoDAOToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.41'
$ws.Range('E48').Value = '  +1.13%  '

$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.16'
$ws.Range('E49').Value = '  -0.06%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '144.89'
$ws.Range('E50').Value = '  +0.19%  '

$ws.Range('E51').Value = '  +0.56%  '
